# Auto-generated Excel COM-interop script
# Applies cached-value updates (market price refresh) to multiple leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1391.5
$ws.Range("I137").Value = 826.26086
$ws.Range("J137").Value = 2156.2354
$ws.Range("K137").Value = 2478.78258
$ws.Range("L137").Value = 6468.706200000001
$ws.Range("M137").Value = 71.21741999999995
$ws.Range("N137").Value = -11568.7062

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 13890336
$ws.Range("I74").Value = 17242668
$ws.Range("J74").Value = 2106.4285
$ws.Range("K74").Value = 17242668
$ws.Range("L74").Value = 2106.4285
$ws.Range("M74").Value = -17241794
$ws.Range("N74").Value = -3854.4285

$ws.Range("H77").Value = 13890336
$ws.Range("I77").Value = 17242668
$ws.Range("J77").Value = 2106.4285
$ws.Range("K77").Value = 86213340
$ws.Range("L77").Value = 10532.1425
$ws.Range("M77").Value = -86208972
$ws.Range("N77").Value = -19268.1425

$ws.Range("H102").Value = 1494.8572
$ws.Range("I102").Value = 1416.4062
$ws.Range("K102").Value = 1416.4062
$ws.Range("M102").Value = 205.5938000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H18").Value = 28330
$ws.Range("J18").Value = 28330
$ws.Range("L18").Value = 28330
$ws.Range("N18").Value = -29388

$ws.Range("H86").Value = 1755.6666
$ws.Range("I86").Value = 1646.8
$ws.Range("J86").Value = 2300
$ws.Range("K86").Value = 1646.8
$ws.Range("L86").Value = 2300
$ws.Range("M86").Value = -523.8
$ws.Range("N86").Value = -4546

$ws.Range("H89").Value = 1755.6666
$ws.Range("I89").Value = 1646.8
$ws.Range("J89").Value = 2300
$ws.Range("K89").Value = 8234
$ws.Range("L89").Value = 11500
$ws.Range("M89").Value = -2618
$ws.Range("N89").Value = -22732

$ws.Range("H105").Value = 1265959
$ws.Range("I105").Value = 2068032.9
$ws.Range("J105").Value = 5557.143
$ws.Range("K105").Value = 2068032.9
$ws.Range("L105").Value = 5557.143
$ws.Range("M105").Value = -2066285.9
$ws.Range("N105").Value = -9051.143

$ws.Range("H107").Value = 683.3571
$ws.Range("I107").Value = 670.56757
$ws.Range("J107").Value = 778
$ws.Range("K107").Value = 670.56757
$ws.Range("L107").Value = 778
$ws.Range("M107").Value = 1249.43243
$ws.Range("N107").Value = -4618

$ws.Range("H134").Value = 3732644
$ws.Range("I134").Value = 5209232.5
$ws.Range("J134").Value = 2315.842
$ws.Range("K134").Value = 15627697.5
$ws.Range("L134").Value = 6947.526
$ws.Range("M134").Value = -15625162.5
$ws.Range("N134").Value = -12017.526

$ws.Range("H135").Value = 42755.645
$ws.Range("J135").Value = 42755.645
$ws.Range("L135").Value = 42755.645
$ws.Range("N135").Value = -52895.645

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3975.8
$ws.Range("I16").Value = 2123.2307
$ws.Range("K16").Value = 2123.2307
$ws.Range("M16").Value = -1836.2307

$ws.Range("H31").Value = 8775052
$ws.Range("I31").Value = 19231844
$ws.Range("J31").Value = 4839.7417
$ws.Range("K31").Value = 19231844
$ws.Range("L31").Value = 4839.7417
$ws.Range("M31").Value = -19231549
$ws.Range("N31").Value = -5429.7417

$ws.Range("H34").Value = 8775052
$ws.Range("I34").Value = 19231844
$ws.Range("J34").Value = 4839.7417
$ws.Range("K34").Value = 19231844
$ws.Range("L34").Value = 4839.7417
$ws.Range("M34").Value = -19231642
$ws.Range("N34").Value = -5243.7417

$ws.Range("H107").Value = 1497.8182
$ws.Range("I107").Value = 673.2308
$ws.Range("J107").Value = 2688.889
$ws.Range("K107").Value = 673.2308
$ws.Range("L107").Value = 2688.889
$ws.Range("M107").Value = 1246.7692
$ws.Range("N107").Value = -6528.889

$ws.Range("H113").Value = 3975.8
$ws.Range("I113").Value = 2123.2307
$ws.Range("K113").Value = 2123.2307
$ws.Range("M113").Value = 46.76929999999993

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 62500510
$ws.Range("I17").Value = 125000024
$ws.Range("J17").Value = 1000
$ws.Range("K17").Value = 375000072
$ws.Range("L17").Value = 3000
$ws.Range("M17").Value = -374999903
$ws.Range("N17").Value = -3338

$ws.Range("H34").Value = 1480.909
$ws.Range("I34").Value = 200
$ws.Range("J34").Value = 1609
$ws.Range("K34").Value = 600
$ws.Range("L34").Value = 4827
$ws.Range("M34").Value = -516
$ws.Range("N34").Value = -4995

$ws.Range("H39").Value = 4026.6667
$ws.Range("I39").Value = 2995
$ws.Range("J39").Value = 4168.9653
$ws.Range("K39").Value = 8985
$ws.Range("L39").Value = 12506.8959
$ws.Range("M39").Value = -8691
$ws.Range("N39").Value = -13094.8959

$ws.Range("H55").Value = 80460.8
$ws.Range("I55").Value = 195027
$ws.Range("J55").Value = 4083.3333
$ws.Range("K55").Value = 585081
$ws.Range("L55").Value = 12249.9999
$ws.Range("M55").Value = -584904
$ws.Range("N55").Value = -12603.9999

$ws.Range("H60").Value = 571.8421
$ws.Range("I60").Value = 263
$ws.Range("K60").Value = 789
$ws.Range("M60").Value = -538

$ws.Range("H123").Value = 1308.5714
$ws.Range("I123").Value = 1193.3334
$ws.Range("J123").Value = 2000
$ws.Range("K123").Value = 3580.0002
$ws.Range("L123").Value = 6000
$ws.Range("M123").Value = -1130.0002
$ws.Range("N123").Value = -10900

$ws.Range("H131").Value = 893.85187
$ws.Range("I131").Value = 342.5
$ws.Range("J131").Value = 989.73914
$ws.Range("K131").Value = 1027.5
$ws.Range("L131").Value = 2969.21742
$ws.Range("M131").Value = 4012.5
$ws.Range("N131").Value = -13049.21742

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = $null
$ws.Range("N19").Value = $null

$ws.Range("H70").Value = 4423.552
$ws.Range("I70").Value = 4331.2383
$ws.Range("K70").Value = 4331.2383
$ws.Range("M70").Value = -4061.2383

$ws.Range("H73").Value = 4423.552
$ws.Range("I73").Value = 4331.2383
$ws.Range("K73").Value = 4331.2383
$ws.Range("M73").Value = -3395.2383

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").Value = $null

$ws.Range("H132").Value = 5420.091
$ws.Range("I132").Value = 5995.9165
$ws.Range("J132").Value = 3884.5557
$ws.Range("K132").Value = 17987.7495
$ws.Range("L132").Value = 11653.6671
$ws.Range("M132").Value = -15457.7495
$ws.Range("N132").Value = -16713.6671

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2143.1353
$ws.Range("I136").Value = 1085.9286
$ws.Range("J136").Value = 5432.222
$ws.Range("K136").Value = 3257.7858
$ws.Range("L136").Value = 16296.666
$ws.Range("M136").Value = -707.7857999999997
$ws.Range("N136").Value = -21396.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 744.44446
$ws.Range("I113").Value = 352.1111
$ws.Range("J113").Value = 940.6111
$ws.Range("K113").Value = 1056.3333
$ws.Range("L113").Value = 2821.8333
$ws.Range("M113").Value = 1113.6667
$ws.Range("N113").Value = -7161.8333

$ws.Range("H132").Value = 1967.4865
$ws.Range("I132").Value = 1478.4445
$ws.Range("J132").Value = 3287.9
$ws.Range("K132").Value = 4435.333500000001
$ws.Range("L132").Value = 9863.7
$ws.Range("M132").Value = -1905.333500000001
$ws.Range("N132").Value = -14923.7
